# Commit: "Modulus opdateret til Modulus Social. Vena har fået status
# godkendt for CPD-DK og XDS Metadata" — but the actual author-visible
# edit that is reflected in the canonical OOXML (per the diff) is a
# simple "new day" rollover of the workbook: the single sheet (and the
# defined name that points at it) is renamed from the 02-12-2025 date
# stamp to 05-12-2025. No cell values/shared strings change in the diff.
#
# Renaming the worksheet via the object model automatically updates the
# sheet-qualified reference inside the workbook-level defined name
# "Lægevagtsystemer" (it tracks the sheet, the same way Excel itself
# keeps defined names in sync when a sheet is renamed through the UI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Opdateret d. 05-12-2025"
